$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix invalid (uppercase) protocol names that were causing an exception
# downstream: MPCX -> mpcx, CPPF -> cppf, GE11 -> ge11
$ws.Range("D5:D44").Value = "mpcx"
$ws.Range("D52:D60").Value = "mpcx"
$ws.Range("D45:D51").Value = "cppf"
$ws.Range("D61:D67").Value = "ge11"

# Update the sheet view: scroll back to the top and select the column
# of protocol names that was just edited
$ws.Range("D5:D67").Select()
